# Slide 5 ("Example") of the deck contains a small process-flow diagram.
# The whole diagram (every shape/connector/picture) is shifted to a new
# position on the slide, and a new caption textbox with a Google Drive
# link is added below it.
#
# NOTE: this COM-interop engine's PowerShell parser does not support
# named ("-Param value") arguments on function/cmdlet calls (they are
# silently dropped / the call becomes a no-op), so every call below uses
# positional arguments only.
#
# EMU -> point conversion used throughout: PowerPoint's Shape.Left/Top/
# Width/Height are expressed in points, while the underlying OOXML
# stores EMUs (1 pt = 12700 EMU).

$EMU_PER_PT = 12700.0

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

function Set-ShapePos([string]$Name, [double]$X, [double]$Y) {
    $sh = $s.Shapes.Item($Name)
    $sh.Left = $X / $EMU_PER_PT
    $sh.Top  = $Y / $EMU_PER_PT
}

function Set-ShapePosExt([string]$Name, [double]$X, [double]$Y, [double]$Cx, [double]$Cy) {
    $sh = $s.Shapes.Item($Name)
    $sh.Left   = $X  / $EMU_PER_PT
    $sh.Top    = $Y  / $EMU_PER_PT
    $sh.Width  = $Cx / $EMU_PER_PT
    $sh.Height = $Cy / $EMU_PER_PT
}

# --- Reposition the existing diagram shapes -------------------------------
Set-ShapePos "Rectangle 3"                    2523744 1693456
Set-ShapePos "TextBox 6"                      2734056 1763988
Set-ShapePos "Rectangle 7"                    874776  2937199
Set-ShapePos "Rectangle 8"                    4163657 2944774
Set-ShapePos "TextBox 9"                      4255007 2944774
Set-ShapePos "TextBox 11"                     1341599 3075311
Set-ShapePos "Straight Arrow Connector 13"    2523744 2133320
Set-ShapePos "Straight Arrow Connector 14"    3791429 2133319
Set-ShapePos "Straight Arrow Connector 18"    2212848 3583529
Set-ShapePos "Rectangle 21"                   874777  4452118
Set-ShapePos "TextBox 22"                     1118702 4613070
Set-ShapePos "Rectangle 23"                   4645064 4452117
Set-ShapePos "TextBox 24"                     4848802 4718016
Set-ShapePos "Straight Arrow Connector 25"    5989232 3572810

# The picture also moves and grows.
Set-ShapePosExt "Picture 27" 7820903 806599 4163833 5276087

# --- Add the new caption textbox with the Google Drive link ---------------
# The source deck's shape-id counter is a simple "next free id" counter
# that is shared per slide and skips ids already present on the slide; to
# reproduce the exact id (32) / default name ("TextBox 31") that PowerPoint
# assigned when the author did this edit, advance the counter with
# throw-away shapes first, then remove them.
$placeholders = @()
for ($i = 0; $i -lt 14; $i++) {
    $placeholders += $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
}
foreach ($ph in $placeholders) {
    $ph.Delete()
}

# Duplicate an existing textbox (rather than Shapes.AddTextbox) so the new
# run inherits the deck's "en-CA" run-formatting context instead of a
# generic "en-US" default.
$source = $s.Shapes.Item("TextBox 11")
$newTextBox = $source.Duplicate().Item(1)
$newTextBox.Name = "TextBox 31"
$newTextBox.TextFrame.TextRange.Text = "https://drive.google.com/file/d/1kj6sNfe_Tib92T1__a-CkNkFFudfbhOt/view?usp=sharing"
$newTextBox.Left   = 163017  / $EMU_PER_PT
$newTextBox.Top    = 6330076 / $EMU_PER_PT
$newTextBox.Width  = 8684514 / $EMU_PER_PT
$newTextBox.Height = 369332  / $EMU_PER_PT

Write-Output ("New shape: id=" + $newTextBox.Id.ToString() + " name=" + $newTextBox.Name)
Write-Output ("Slide 5 shape count=" + $s.Shapes.Count.ToString())
